# "Update decoding and data"
#
# Clears the LSB-decode "marker" value ("x") from column E (the "Sequential"
# decode-order column) for every image row whose D/E/F triple had all three
# markers set (rows 5, 17, 21, 25, 29, 31, 32, 33, 34 on the "Images" sheet),
# then moves the active selection on that sheet to E30 and nudges the
# workbook window geometry.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Images")

$rows = @(5, 17, 21, 25, 29, 31, 32, 33, 34)
foreach ($r in $rows) {
    $ws.Range("E$r").Value = ""
}

# Move the sheet's active selection (was L9) to E30.
$ws.Activate()
$ws.Range("E30").Select()

# Match the updated workbook window geometry.
$win = $excel.ActiveWindow
$win.Left = 980
$win.Top = 500
$win.Width = 18700
$win.Height = 21100
